$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column C: copy header cell formatting from B1, then set the new date ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Cells.Item(1, 3).Value = "13-01-2023"

# --- Fund rows re-ordered alphabetically (rows 2-40), with "avg" (41) and "total" (42) moved to the bottom ---
# Column A: fund/label name, Column B: existing balance (06-01-2023), Column C: new balance (13-01-2023)

$ws.Cells.Item(2, 1).Value = "1810 Renta variable"
$ws.Cells.Item(2, 2).Value = 94509.85
$ws.Cells.Item(2, 3).Value = 94534.46

$ws.Cells.Item(3, 1).Value = "1822 Raices Valores Negociables"
$ws.Cells.Item(3, 2).Value = 238818.2
$ws.Cells.Item(3, 3).Value = 239270.73

$ws.Cells.Item(4, 1).Value = "Adcap IOL Acciones Argentina"
$ws.Cells.Item(4, 2).Value = 29313.29
$ws.Cells.Item(4, 3).Value = 31308.64

$ws.Cells.Item(5, 1).Value = "Allaria Acciones"
$ws.Cells.Item(5, 2).Value = 17364.12
$ws.Cells.Item(5, 3).Value = 17412.14

$ws.Cells.Item(6, 1).Value = "Alpha Acciones"
$ws.Cells.Item(6, 2).Value = 79846.41
$ws.Cells.Item(6, 3).Value = 98190.67

$ws.Cells.Item(7, 1).Value = "Alpha Mega"
$ws.Cells.Item(7, 2).Value = 328410.27
$ws.Cells.Item(7, 3).Value = 328385.78

$ws.Cells.Item(8, 1).Value = "Alpha Mercosur"
$ws.Cells.Item(8, 2).Value = 579063.9
$ws.Cells.Item(8, 3).Value = 609869.29

$ws.Cells.Item(9, 1).Value = "Alpha planeam equil"
$ws.Cells.Item(9, 2).Value = 11574.97
$ws.Cells.Item(9, 3).Value = 7285.59

$ws.Cells.Item(10, 1).Value = "Alpha renta balan global"
$ws.Cells.Item(10, 2).Value = 646692.49
$ws.Cells.Item(10, 3).Value = 646095.37

$ws.Cells.Item(11, 1).Value = "Argenfunds"
$ws.Cells.Item(11, 2).Value = 9504.72
$ws.Cells.Item(11, 3).Value = 9488.53

$ws.Cells.Item(12, 1).Value = "Arpenta ex Mercosur"
$ws.Cells.Item(12, 2).Value = 19996.28
$ws.Cells.Item(12, 3).Value = 20004.2

$ws.Cells.Item(13, 1).Value = "Balanz"
$ws.Cells.Item(13, 2).Value = 146514.67
$ws.Cells.Item(13, 3).Value = 233548.68

$ws.Cells.Item(14, 1).Value = "CMA acciones"
$ws.Cells.Item(14, 2).Value = 187839.74
$ws.Cells.Item(14, 3).Value = 187844.01

$ws.Cells.Item(15, 1).Value = "Compass Crecimiento"
$ws.Cells.Item(15, 2).Value = 731281.3
$ws.Cells.Item(15, 3).Value = 729452.82

$ws.Cells.Item(16, 1).Value = "Consultatio Renta Variable"
$ws.Cells.Item(16, 2).Value = 10230.37
$ws.Cells.Item(16, 3).Value = 9913.58

$ws.Cells.Item(17, 1).Value = "Delta Acciones"
$ws.Cells.Item(17, 2).Value = 42544.86
$ws.Cells.Item(17, 3).Value = 42414.27

$ws.Cells.Item(18, 1).Value = "Delta Select"
$ws.Cells.Item(18, 2).Value = 277593.58
$ws.Cells.Item(18, 3).Value = 277861.58

$ws.Cells.Item(19, 1).Value = "Delta gestion V"
$ws.Cells.Item(19, 2).Value = 25359.32
$ws.Cells.Item(19, 3).Value = 25472.07

$ws.Cells.Item(20, 1).Value = "Fima Acciones"
$ws.Cells.Item(20, 2).Value = 184135.65
$ws.Cells.Item(20, 3).Value = 183557.94

$ws.Cells.Item(21, 1).Value = "Fima PB Acciones"
$ws.Cells.Item(21, 2).Value = 219846.69
$ws.Cells.Item(21, 3).Value = 230285.08

$ws.Cells.Item(22, 1).Value = "Gainvest Renta Variable"
$ws.Cells.Item(22, 2).Value = 54382.44
$ws.Cells.Item(22, 3).Value = 54419.07

$ws.Cells.Item(23, 1).Value = "Goal Acciones Argentinas"
$ws.Cells.Item(23, 2).Value = 23186.18
$ws.Cells.Item(23, 3).Value = 23208.67

$ws.Cells.Item(24, 1).Value = "Goal acciones plus"
$ws.Cells.Item(24, 2).Value = 5375.02
$ws.Cells.Item(24, 3).Value = 5374.26

$ws.Cells.Item(25, 1).Value = "HF Acciones Argentinas"
$ws.Cells.Item(25, 2).Value = 95060.08
$ws.Cells.Item(25, 3).Value = 95107.04

$ws.Cells.Item(26, 1).Value = "HF Acciones Lideres"
$ws.Cells.Item(26, 2).Value = 102176.73
$ws.Cells.Item(26, 3).Value = 112199.76

$ws.Cells.Item(27, 1).Value = "IAM Renta Variable"
$ws.Cells.Item(27, 2).Value = 30650.72
$ws.Cells.Item(27, 3).Value = 32094.72

$ws.Cells.Item(28, 1).Value = "IEB Value"
$ws.Cells.Item(28, 2).Value = 3618.11
$ws.Cells.Item(28, 3).Value = 3617.28

$ws.Cells.Item(29, 1).Value = "Lombardi"
$ws.Cells.Item(29, 2).Value = 37116.79
$ws.Cells.Item(29, 3).Value = 41361.28

$ws.Cells.Item(30, 1).Value = "MAF"
$ws.Cells.Item(30, 2).Value = 16184.43
$ws.Cells.Item(30, 3).Value = 16218.33

$ws.Cells.Item(31, 1).Value = "Megainver"
$ws.Cells.Item(31, 2).Value = 26248.2
$ws.Cells.Item(31, 3).Value = 26182.15

$ws.Cells.Item(32, 1).Value = "Pellegrini Acciones"
$ws.Cells.Item(32, 2).Value = 69990.33
$ws.Cells.Item(32, 3).Value = 79872.56

$ws.Cells.Item(33, 1).Value = "Pionero Acciones"
$ws.Cells.Item(33, 2).Value = 98147.54
$ws.Cells.Item(33, 3).Value = 98071.52

$ws.Cells.Item(34, 1).Value = "Premier Renta Variable"
$ws.Cells.Item(34, 2).Value = 58260.09
$ws.Cells.Item(34, 3).Value = 58343.03

$ws.Cells.Item(35, 1).Value = "Quinquela Acciones"
$ws.Cells.Item(35, 2).Value = 84985.54
$ws.Cells.Item(35, 3).Value = 84911.81

$ws.Cells.Item(36, 1).Value = "Rofex 20 Renta Variable"
$ws.Cells.Item(36, 2).Value = 60333.63
$ws.Cells.Item(36, 3).Value = 60354.15

$ws.Cells.Item(37, 1).Value = "Supefondo RV"
$ws.Cells.Item(37, 2).Value = 1120946.22
$ws.Cells.Item(37, 3).Value = 1119813.64

$ws.Cells.Item(38, 1).Value = "Superfondo "
$ws.Cells.Item(38, 2).Value = 1181522.05
$ws.Cells.Item(38, 3).Value = 1182138.57

$ws.Cells.Item(39, 1).Value = "Supergestion"
$ws.Cells.Item(39, 2).Value = 469846.65
$ws.Cells.Item(39, 3).Value = 470155.03

$ws.Cells.Item(40, 1).Value = "Toronto Trust Multimercado"
$ws.Cells.Item(40, 2).Value = 37010.05
$ws.Cells.Item(40, 3).Value = 37052.93

$ws.Cells.Item(41, 1).Value = "avg"
$ws.Cells.Item(41, 2).Value = 191166.19
$ws.Cells.Item(41, 3).Value = 195453.62

$ws.Cells.Item(42, 1).Value = "total"
$ws.Cells.Item(42, 2).Value = 7455481.48
$ws.Cells.Item(42, 3).Value = 7622691.23
